$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.040.84'
$ws.Range('E2').Value = '  -0.37%  '

$ws.Range('D3').Value = '2.213.56'
$ws.Range('E3').Value = '  -1.20%  '

$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.77'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.23%  '

$ws.Range('E6').Value = '  -0.54%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.08'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.06%  '

$ws.Range('E8').Value = '  +0.12%  '

$ws.Range('E9').Value = '  -2.09%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.23'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.13%  '

$ws.Range('E11').Value = '  +1.07%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.05'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.15%  '

$ws.Range('E13').Value = '  -0.32%  '

$ws.Range('D14').Value = '2.545.28'
$ws.Range('E14').Value = '  -1.15%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.16'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.20%  '

$ws.Range('E16').Value = '  -1.78%  '

$ws.Range('D17').Value = '2.274.15'
$ws.Range('E17').Value = '  +1.65%  '

$ws.Range('D18').Value = '41.875.09'
$ws.Range('E18').Value = '  -0.39%  '

$ws.Range('E19').Value = '  +9.22%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.73'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.97%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.11'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.48%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.19'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +15.46%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.22'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.89%  '

$ws.Range('E24').Value = '  -6.66%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.69'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.88%  '

$ws.Range('E26').Value = '  +0.14%  '

$ws.Range('E27').Value = '  -0.79%  '

$ws.Range('E28').Value = '  -2.10%  '

$ws.Range('E29').Value = '  +4.59%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.01'
$ws.Range('D30').ClearFormats()

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.56'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.26%  '

$ws.Range('E32').Value = '  +7.97%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0787'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.62%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '28.69'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.40%  '

$ws.Range('E36').Value = '  -7.97%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.24'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.23%  '

$ws.Range('E38').Value = '  -4.62%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.23'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.12%  '

$ws.Range('B40').Value = 'MultiversX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '65.15'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.98%  '

$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.12'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.99%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.61'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.57%  '

$ws.Range('E43').Value = '  -3.35%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.68'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.52%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.78'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.68%  '

$ws.Range('E46').Value = '  -2.05%  '

$ws.Range('E47').Value = '  +5.10%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.11'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.62%  '

$ws.Range('E49').Value = '  -0.14%  '

$ws.Range('E50').Value = '  -0.18%  '

$ws.Range('D51').Value = '2.418.07'
